$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 873.6
$ws.Cells.Item(28, 9).Value = 1056.3334
$ws.Cells.Item(28, 11).Value = 1056.3334
$ws.Cells.Item(28, 13).Value = -571.3334
$ws.Cells.Item(32, 8).Value = 1500
$ws.Cells.Item(32, 9).Value = 2000
$ws.Cells.Item(32, 10).Value = 1333.3334
$ws.Cells.Item(32, 11).Value = 2000
$ws.Cells.Item(32, 12).Value = 1333.3334
$ws.Cells.Item(32, 13).Value = -1674
$ws.Cells.Item(32, 14).Value = -1985.3334
$ws.Cells.Item(33, 8).Value = 437.66666
$ws.Cells.Item(33, 9).Value = 415.3846
$ws.Cells.Item(33, 10).Value = 582.5
$ws.Cells.Item(33, 11).Value = 415.3846
$ws.Cells.Item(33, 12).Value = 582.5
$ws.Cells.Item(33, 13).Value = -186.3846
$ws.Cells.Item(33, 14).Value = -1040.5
$ws.Cells.Item(38, 8).Value = 29.88889
$ws.Cells.Item(38, 9).Value = 27.375
$ws.Cells.Item(38, 10).Value = 50
$ws.Cells.Item(38, 11).Value = 82.125
$ws.Cells.Item(38, 12).Value = 150
$ws.Cells.Item(38, 13).Value = 289.875
$ws.Cells.Item(38, 14).Value = -894
$ws.Cells.Item(61, 8).Value = 450.4
$ws.Cells.Item(61, 9).Value = 450.4
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 1351.2
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -1179.2
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(62, 8).Value = 7358.75
$ws.Cells.Item(62, 9).Value = 4997
$ws.Cells.Item(62, 11).Value = 4997
$ws.Cells.Item(62, 13).Value = -4373
$ws.Cells.Item(65, 8).Value = 7358.75
$ws.Cells.Item(65, 9).Value = 4997
$ws.Cells.Item(65, 11).Value = 24985
$ws.Cells.Item(65, 13).Value = -21865
$ws.Cells.Item(70, 8).Value = 3454.5454
$ws.Cells.Item(70, 10).Value = 3000
$ws.Cells.Item(70, 12).Value = 9000
$ws.Cells.Item(70, 14).Value = -9540
$ws.Cells.Item(73, 8).Value = 3454.5454
$ws.Cells.Item(73, 10).Value = 3000
$ws.Cells.Item(73, 12).Value = 9000
$ws.Cells.Item(73, 14).Value = -10872
$ws.Cells.Item(86, 8).Value = 2359.2104
$ws.Cells.Item(86, 10).Value = 1543.5
$ws.Cells.Item(86, 12).Value = 1543.5
$ws.Cells.Item(86, 14).Value = -3789.5
$ws.Cells.Item(88, 8).Value = 3824.1
$ws.Cells.Item(88, 9).Value = 3622.5
$ws.Cells.Item(88, 10).Value = 3874.5
$ws.Cells.Item(88, 11).Value = 3622.5
$ws.Cells.Item(88, 12).Value = 3874.5
$ws.Cells.Item(88, 13).Value = -3216.5
$ws.Cells.Item(88, 14).Value = -4686.5
$ws.Cells.Item(89, 8).Value = 2359.2104
$ws.Cells.Item(89, 10).Value = 1543.5
$ws.Cells.Item(89, 12).Value = 7717.5
$ws.Cells.Item(89, 14).Value = -18949.5
$ws.Cells.Item(91, 8).Value = 3824.1
$ws.Cells.Item(91, 9).Value = 3622.5
$ws.Cells.Item(91, 10).Value = 3874.5
$ws.Cells.Item(91, 11).Value = 3622.5
$ws.Cells.Item(91, 12).Value = 3874.5
$ws.Cells.Item(91, 13).Value = -2218.5
$ws.Cells.Item(91, 14).Value = -6682.5
$ws.Cells.Item(92, 8).Value = 100.545456
$ws.Cells.Item(92, 9).Value = 81.625
$ws.Cells.Item(92, 11).Value = 81.625
$ws.Cells.Item(92, 13).Value = 1166.375
$ws.Cells.Item(94, 8).Value = 5862.5557
$ws.Cells.Item(94, 9).Value = 2220.375
$ws.Cells.Item(94, 10).Value = 35000
$ws.Cells.Item(94, 11).Value = 2220.375
$ws.Cells.Item(94, 12).Value = 35000
$ws.Cells.Item(94, 13).Value = -1769.375
$ws.Cells.Item(94, 14).Value = -35902
$ws.Cells.Item(97, 8).Value = 4949.5
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 4949.5
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 14848.5
$ws.Cells.Item(97, 13).ClearContents()
$ws.Cells.Item(97, 14).Value = -15840.5
$ws.Cells.Item(100, 8).Value = 2046.7646
$ws.Cells.Item(100, 9).Value = 2190.0667
$ws.Cells.Item(100, 10).Value = 972
$ws.Cells.Item(100, 11).Value = 2190.0667
$ws.Cells.Item(100, 12).Value = 972
$ws.Cells.Item(100, 13).Value = -1649.0667
$ws.Cells.Item(100, 14).Value = -2054
$ws.Cells.Item(103, 8).Value = 2071.3333
$ws.Cells.Item(103, 10).Value = 2551.818
$ws.Cells.Item(103, 12).Value = 7655.454000000001
$ws.Cells.Item(103, 14).Value = -8827.454000000002
$ws.Cells.Item(106, 8).Value = 1498
$ws.Cells.Item(106, 9).Value = 1498
$ws.Cells.Item(106, 11).Value = 1498
$ws.Cells.Item(106, 13).Value = -867
$ws.Cells.Item(113, 8).Value = 3499.5
$ws.Cells.Item(113, 9).Value = 3499.5
$ws.Cells.Item(113, 11).Value = 3499.5
$ws.Cells.Item(113, 13).Value = -245.5
$ws.Cells.Item(132, 8).Value = 3882.8386
$ws.Cells.Item(132, 9).Value = 3215.8
$ws.Cells.Item(132, 10).Value = 5095.636
$ws.Cells.Item(132, 11).Value = 9647.400000000001
$ws.Cells.Item(132, 12).Value = 15286.908
$ws.Cells.Item(132, 13).Value = -7117.400000000001
$ws.Cells.Item(132, 14).Value = -20346.908
$ws.Cells.Item(135, 8).Value = 298.25
$ws.Cells.Item(135, 9).Value = 298.25
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 2684.25
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -149.25
$ws.Cells.Item(135, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 2629.4707
$ws.Cells.Item(138, 9).Value = 2446.8
$ws.Cells.Item(138, 11).Value = 7340.400000000001
$ws.Cells.Item(138, 13).Value = -2200.400000000001
$ws.Cells.Item(141, 8).Value = 1454.9286
$ws.Cells.Item(141, 9).Value = 1489.9231
$ws.Cells.Item(141, 11).Value = 4469.7693
$ws.Cells.Item(141, 13).Value = 710.2307000000001

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 9803.6
$ws.Cells.Item(2, 10).Value = 10004.5
$ws.Cells.Item(2, 12).Value = 10004.5
$ws.Cells.Item(2, 14).Value = -10230.5
$ws.Cells.Item(27, 8).Value = 11496
$ws.Cells.Item(27, 10).Value = 11496
$ws.Cells.Item(27, 12).Value = 11496
$ws.Cells.Item(27, 14).Value = -11864
$ws.Cells.Item(32, 8).Value = 3299.353
$ws.Cells.Item(32, 9).Value = 2239.3667
$ws.Cells.Item(32, 11).Value = 2239.3667
$ws.Cells.Item(32, 13).Value = -1952.3667
$ws.Cells.Item(45, 8).Value = 3773.7334
$ws.Cells.Item(45, 9).Value = 2332.3333
$ws.Cells.Item(45, 10).Value = 4734.6665
$ws.Cells.Item(45, 11).Value = 2332.3333
$ws.Cells.Item(45, 12).Value = 4734.6665
$ws.Cells.Item(45, 13).Value = -1955.3333
$ws.Cells.Item(45, 14).Value = -5488.6665
$ws.Cells.Item(50, 8).Value = 22625
$ws.Cells.Item(50, 9).Value = 20000
$ws.Cells.Item(50, 10).Value = 25250
$ws.Cells.Item(50, 11).Value = 20000
$ws.Cells.Item(50, 12).Value = 25250
$ws.Cells.Item(50, 13).Value = -19286
$ws.Cells.Item(50, 14).Value = -26678
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 14).ClearContents()
$ws.Cells.Item(55, 8).Value = 25000
$ws.Cells.Item(55, 10).Value = 25000
$ws.Cells.Item(55, 12).Value = 25000
$ws.Cells.Item(55, 14).Value = -25630
$ws.Cells.Item(61, 8).Value = 986.1429000000001
$ws.Cells.Item(61, 9).Value = 986.1429000000001
$ws.Cells.Item(61, 11).Value = 986.1429000000001
$ws.Cells.Item(61, 13).Value = -774.1429000000001
$ws.Cells.Item(74, 8).Value = 1430.6
$ws.Cells.Item(74, 9).Value = 1466.3077
$ws.Cells.Item(74, 10).Value = 1198.5
$ws.Cells.Item(74, 11).Value = 1466.3077
$ws.Cells.Item(74, 12).Value = 1198.5
$ws.Cells.Item(74, 13).Value = -592.3077000000001
$ws.Cells.Item(74, 14).Value = -2946.5
$ws.Cells.Item(77, 8).Value = 1430.6
$ws.Cells.Item(77, 9).Value = 1466.3077
$ws.Cells.Item(77, 10).Value = 1198.5
$ws.Cells.Item(77, 11).Value = 7331.538500000001
$ws.Cells.Item(77, 12).Value = 5992.5
$ws.Cells.Item(77, 13).Value = -2963.538500000001
$ws.Cells.Item(77, 14).Value = -14728.5
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(97, 8).Value = 797.9
$ws.Cells.Item(97, 9).Value = 747.5294
$ws.Cells.Item(97, 10).Value = 1083.3334
$ws.Cells.Item(97, 11).Value = 747.5294
$ws.Cells.Item(97, 12).Value = 1083.3334
$ws.Cells.Item(97, 13).Value = -251.5294
$ws.Cells.Item(97, 14).Value = -2075.3334
$ws.Cells.Item(102, 8).Value = 1854.5
$ws.Cells.Item(102, 9).Value = 1119.4286
$ws.Cells.Item(102, 11).Value = 1119.4286
$ws.Cells.Item(102, 13).Value = 502.5714
$ws.Cells.Item(110, 8).Value = 968.55554
$ws.Cells.Item(110, 9).Value = 968.55554
$ws.Cells.Item(110, 11).Value = 968.55554
$ws.Cells.Item(110, 13).Value = 1076.44446
$ws.Cells.Item(116, 8).Value = 9803.6
$ws.Cells.Item(116, 10).Value = 10004.5
$ws.Cells.Item(116, 12).Value = 10004.5
$ws.Cells.Item(116, 14).Value = -14592.5
$ws.Cells.Item(122, 8).Value = 4956
$ws.Cells.Item(122, 9).Value = 4956
$ws.Cells.Item(122, 11).Value = 14868
$ws.Cells.Item(122, 13).Value = -12418
$ws.Cells.Item(132, 8).Value = 1081
$ws.Cells.Item(132, 9).Value = 941.3333
$ws.Cells.Item(132, 10).Value = 1500
$ws.Cells.Item(132, 11).Value = 2823.9999
$ws.Cells.Item(132, 12).Value = 4500
$ws.Cells.Item(132, 13).Value = -293.9998999999998
$ws.Cells.Item(132, 14).Value = -9560
$ws.Cells.Item(136, 8).Value = 986.1429000000001
$ws.Cells.Item(136, 9).Value = 986.1429000000001
$ws.Cells.Item(136, 11).Value = 2958.4287
$ws.Cells.Item(136, 13).Value = -408.4287000000004

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 9803.6
$ws.Cells.Item(3, 10).Value = 10004.5
$ws.Cells.Item(3, 12).Value = 10004.5
$ws.Cells.Item(3, 14).Value = -10232.5
$ws.Cells.Item(7, 8).Value = 3800269.5
$ws.Cells.Item(7, 9).Value = 3800269.5
$ws.Cells.Item(7, 11).Value = 3800269.5
$ws.Cells.Item(7, 13).Value = -3800156.5
$ws.Cells.Item(20, 8).Value = 1259.5834
$ws.Cells.Item(20, 9).Value = 1311.6
$ws.Cells.Item(20, 11).Value = 1311.6
$ws.Cells.Item(20, 13).Value = -1064.6
$ws.Cells.Item(22, 8).Value = 970
$ws.Cells.Item(22, 9).Value = 970
$ws.Cells.Item(22, 11).Value = 970
$ws.Cells.Item(22, 13).Value = -797
$ws.Cells.Item(51, 8).Value = 99998.5
$ws.Cells.Item(51, 10).Value = 99998.5
$ws.Cells.Item(51, 12).Value = 99998.5
$ws.Cells.Item(51, 14).Value = -100980.5
$ws.Cells.Item(54, 8).Value = 9065.833000000001
$ws.Cells.Item(54, 9).Value = 7850
$ws.Cells.Item(54, 10).Value = 11497.5
$ws.Cells.Item(54, 11).Value = 7850
$ws.Cells.Item(54, 12).Value = 11497.5
$ws.Cells.Item(54, 13).Value = -7366
$ws.Cells.Item(54, 14).Value = -12465.5
$ws.Cells.Item(64, 8).Value = 966.3333
$ws.Cells.Item(64, 9).Value = 903.3333
$ws.Cells.Item(64, 11).Value = 903.3333
$ws.Cells.Item(64, 13).Value = -678.3333
$ws.Cells.Item(67, 8).Value = 966.3333
$ws.Cells.Item(67, 9).Value = 903.3333
$ws.Cells.Item(67, 11).Value = 903.3333
$ws.Cells.Item(67, 13).Value = -123.3333
$ws.Cells.Item(86, 8).Value = 3933.7727
$ws.Cells.Item(86, 9).Value = 3306.2666
$ws.Cells.Item(86, 10).Value = 5278.4287
$ws.Cells.Item(86, 11).Value = 3306.2666
$ws.Cells.Item(86, 12).Value = 5278.4287
$ws.Cells.Item(86, 13).Value = -2183.2666
$ws.Cells.Item(86, 14).Value = -7524.4287
$ws.Cells.Item(89, 8).Value = 3933.7727
$ws.Cells.Item(89, 9).Value = 3306.2666
$ws.Cells.Item(89, 10).Value = 5278.4287
$ws.Cells.Item(89, 11).Value = 16531.333
$ws.Cells.Item(89, 12).Value = 26392.1435
$ws.Cells.Item(89, 13).Value = -10915.333
$ws.Cells.Item(89, 14).Value = -37624.14350000001
$ws.Cells.Item(94, 8).Value = 3245.05
$ws.Cells.Item(94, 9).Value = 2660.0667
$ws.Cells.Item(94, 11).Value = 2660.0667
$ws.Cells.Item(94, 13).Value = -2209.0667
$ws.Cells.Item(95, 8).Value = 49312
$ws.Cells.Item(95, 10).Value = 49312
$ws.Cells.Item(95, 12).Value = 49312
$ws.Cells.Item(95, 14).Value = -54804
$ws.Cells.Item(99, 8).Value = 5862.5713
$ws.Cells.Item(99, 9).Value = 4859.75
$ws.Cells.Item(99, 10).Value = 7199.6665
$ws.Cells.Item(99, 11).Value = 4859.75
$ws.Cells.Item(99, 12).Value = 7199.6665
$ws.Cells.Item(99, 13).Value = -3361.75
$ws.Cells.Item(99, 14).Value = -10195.6665
$ws.Cells.Item(105, 8).Value = 3168.625
$ws.Cells.Item(105, 9).Value = 2308.1667
$ws.Cells.Item(105, 11).Value = 2308.1667
$ws.Cells.Item(105, 13).Value = -561.1667000000002
$ws.Cells.Item(134, 8).Value = 3233.6924
$ws.Cells.Item(134, 9).Value = 3276.2727
$ws.Cells.Item(134, 11).Value = 9828.8181
$ws.Cells.Item(134, 13).Value = -7293.8181

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1673.6666
$ws.Cells.Item(31, 9).Value = 1270.8
$ws.Cells.Item(31, 10).Value = 2479.4
$ws.Cells.Item(31, 11).Value = 1270.8
$ws.Cells.Item(31, 12).Value = 2479.4
$ws.Cells.Item(31, 13).Value = -975.8
$ws.Cells.Item(31, 14).Value = -3069.4
$ws.Cells.Item(34, 8).Value = 1673.6666
$ws.Cells.Item(34, 9).Value = 1270.8
$ws.Cells.Item(34, 10).Value = 2479.4
$ws.Cells.Item(34, 11).Value = 1270.8
$ws.Cells.Item(34, 12).Value = 2479.4
$ws.Cells.Item(34, 13).Value = -1068.8
$ws.Cells.Item(34, 14).Value = -2883.4
$ws.Cells.Item(35, 8).Value = 12920.667
$ws.Cells.Item(35, 9).Value = 11631
$ws.Cells.Item(35, 10).Value = 15500
$ws.Cells.Item(35, 11).Value = 11631
$ws.Cells.Item(35, 12).Value = 15500
$ws.Cells.Item(35, 13).Value = -11337
$ws.Cells.Item(35, 14).Value = -16088
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 14).ClearContents()
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 14).ClearContents()
$ws.Cells.Item(50, 8).Value = 34797
$ws.Cells.Item(50, 10).Value = 39995
$ws.Cells.Item(50, 12).Value = 39995
$ws.Cells.Item(50, 14).Value = -41245
$ws.Cells.Item(55, 8).Value = 6000
$ws.Cells.Item(55, 9).Value = 6000
$ws.Cells.Item(55, 11).Value = 6000
$ws.Cells.Item(55, 13).Value = -5685
$ws.Cells.Item(58, 8).Value = 1404.6875
$ws.Cells.Item(58, 9).Value = 1185
$ws.Cells.Item(58, 11).Value = 1185
$ws.Cells.Item(58, 13).Value = -982
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()
$ws.Cells.Item(86, 9).Value = 34851484
$ws.Cells.Item(86, 11).Value = 34851484
$ws.Cells.Item(86, 13).Value = -34850361
$ws.Cells.Item(89, 9).Value = 34851484
$ws.Cells.Item(89, 11).Value = 174257420
$ws.Cells.Item(89, 13).Value = -174251804
$ws.Cells.Item(99, 8).Value = 4405.2
$ws.Cells.Item(99, 9).Value = 5071.364
$ws.Cells.Item(99, 10).Value = 3591
$ws.Cells.Item(99, 11).Value = 5071.364
$ws.Cells.Item(99, 12).Value = 3591
$ws.Cells.Item(99, 13).Value = -3573.364
$ws.Cells.Item(99, 14).Value = -6587
$ws.Cells.Item(105, 8).Value = 2600.8
$ws.Cells.Item(105, 9).Value = 2752.5
$ws.Cells.Item(105, 10).Value = 2499.6667
$ws.Cells.Item(105, 11).Value = 2752.5
$ws.Cells.Item(105, 12).Value = 2499.6667
$ws.Cells.Item(105, 13).Value = -1005.5
$ws.Cells.Item(105, 14).Value = -5993.6667
$ws.Cells.Item(122, 8).Value = 2252
$ws.Cells.Item(122, 9).Value = 2252
$ws.Cells.Item(122, 11).Value = 6756
$ws.Cells.Item(122, 13).Value = -4306
$ws.Cells.Item(126, 8).Value = 4405.2
$ws.Cells.Item(126, 9).Value = 5071.364
$ws.Cells.Item(126, 10).Value = 3591
$ws.Cells.Item(126, 11).Value = 15214.092
$ws.Cells.Item(126, 12).Value = 10773
$ws.Cells.Item(126, 13).Value = -12744.092
$ws.Cells.Item(126, 14).Value = -15713
$ws.Cells.Item(134, 8).Value = 3497.5715
$ws.Cells.Item(134, 9).Value = 2896.6
$ws.Cells.Item(134, 11).Value = 8689.799999999999
$ws.Cells.Item(134, 13).Value = -6154.799999999999
$ws.Cells.Item(136, 8).Value = 1404.6875
$ws.Cells.Item(136, 9).Value = 1185
$ws.Cells.Item(136, 11).Value = 3555
$ws.Cells.Item(136, 13).Value = -1005

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 391.5
$ws.Cells.Item(7, 9).Value = 87
$ws.Cells.Item(7, 10).Value = 696
$ws.Cells.Item(7, 11).Value = 261
$ws.Cells.Item(7, 12).Value = 2088
$ws.Cells.Item(7, 13).Value = -149
$ws.Cells.Item(7, 14).Value = -2312
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(68, 14).ClearContents()
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).ClearContents()
$ws.Cells.Item(71, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 350
$ws.Cells.Item(113, 9).Value = 300
$ws.Cells.Item(113, 11).Value = 900
$ws.Cells.Item(113, 13).Value = 1270
$ws.Cells.Item(122, 8).Value = 491.2143
$ws.Cells.Item(122, 9).Value = 512.1667
$ws.Cells.Item(122, 10).Value = 475.5
$ws.Cells.Item(122, 11).Value = 4609.5003
$ws.Cells.Item(122, 12).Value = 4279.5
$ws.Cells.Item(122, 13).Value = -2159.5003
$ws.Cells.Item(122, 14).Value = -9179.5
$ws.Cells.Item(123, 8).Value = 1999
$ws.Cells.Item(123, 9).Value = 1999
$ws.Cells.Item(123, 11).Value = 5997
$ws.Cells.Item(123, 13).Value = -3547
$ws.Cells.Item(128, 8).Value = 278575.66
$ws.Cells.Item(128, 9).Value = 278575.66
$ws.Cells.Item(128, 11).Value = 835726.98
$ws.Cells.Item(128, 13).Value = -830746.98
$ws.Cells.Item(132, 8).Value = 458.33334
$ws.Cells.Item(132, 9).Value = 450
$ws.Cells.Item(132, 10).Value = 462.5
$ws.Cells.Item(132, 11).Value = 4050
$ws.Cells.Item(132, 12).Value = 4162.5
$ws.Cells.Item(132, 13).Value = -1520
$ws.Cells.Item(132, 14).Value = -9222.5

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 500
$ws.Cells.Item(5, 9).Value = 500
$ws.Cells.Item(5, 11).Value = 500
$ws.Cells.Item(5, 13).Value = -388
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(54, 8).Value = 14994
$ws.Cells.Item(54, 10).Value = 14994
$ws.Cells.Item(54, 12).Value = 14994
$ws.Cells.Item(54, 14).Value = -15774
$ws.Cells.Item(70, 8).Value = 4995.5
$ws.Cells.Item(70, 9).Value = 4996
$ws.Cells.Item(70, 10).Value = 4994
$ws.Cells.Item(70, 11).Value = 4996
$ws.Cells.Item(70, 12).Value = 4994
$ws.Cells.Item(70, 13).Value = -4726
$ws.Cells.Item(70, 14).Value = -5534
$ws.Cells.Item(73, 8).Value = 4995.5
$ws.Cells.Item(73, 9).Value = 4996
$ws.Cells.Item(73, 10).Value = 4994
$ws.Cells.Item(73, 11).Value = 4996
$ws.Cells.Item(73, 12).Value = 4994
$ws.Cells.Item(73, 13).Value = -4060
$ws.Cells.Item(73, 14).Value = -6866
$ws.Cells.Item(80, 8).Value = 2812.5
$ws.Cells.Item(80, 9).Value = 2649.5
$ws.Cells.Item(80, 10).Value = 2975.5
$ws.Cells.Item(80, 11).Value = 2649.5
$ws.Cells.Item(80, 12).Value = 2975.5
$ws.Cells.Item(80, 13).Value = -1651.5
$ws.Cells.Item(80, 14).Value = -4971.5
$ws.Cells.Item(83, 8).Value = 2812.5
$ws.Cells.Item(83, 9).Value = 2649.5
$ws.Cells.Item(83, 10).Value = 2975.5
$ws.Cells.Item(83, 11).Value = 13247.5
$ws.Cells.Item(83, 12).Value = 14877.5
$ws.Cells.Item(83, 13).Value = -8255.5
$ws.Cells.Item(83, 14).Value = -24861.5
$ws.Cells.Item(97, 8).Value = 858.2143
$ws.Cells.Item(97, 9).Value = 575.2727
$ws.Cells.Item(97, 11).Value = 575.2727
$ws.Cells.Item(97, 13).Value = -79.27269999999999
$ws.Cells.Item(102, 8).Value = 920.8333
$ws.Cells.Item(102, 9).Value = 723.5714
$ws.Cells.Item(102, 11).Value = 723.5714
$ws.Cells.Item(102, 13).Value = 898.4286
$ws.Cells.Item(107, 8).Value = 365.33334
$ws.Cells.Item(107, 9).Value = 307.72726
$ws.Cells.Item(107, 11).Value = 307.72726
$ws.Cells.Item(107, 13).Value = 1612.27274
$ws.Cells.Item(113, 8).Value = 2175.4
$ws.Cells.Item(113, 9).Value = 2175.4
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2175.4
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -5.400000000000091
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 13).Value = -3550
$ws.Cells.Item(126, 8).Value = 2599.875
$ws.Cells.Item(126, 9).Value = 2340.2
$ws.Cells.Item(126, 10).Value = 3032.6667
$ws.Cells.Item(126, 11).Value = 7020.599999999999
$ws.Cells.Item(126, 12).Value = 9098.000100000001
$ws.Cells.Item(126, 13).Value = -4550.599999999999
$ws.Cells.Item(126, 14).Value = -14038.0001
$ws.Cells.Item(132, 8).Value = 4811.4287
$ws.Cells.Item(132, 9).Value = 4336
$ws.Cells.Item(132, 11).Value = 13008
$ws.Cells.Item(132, 13).Value = -10478

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 850
$ws.Cells.Item(2, 9).Value = 850
$ws.Cells.Item(2, 11).Value = 850
$ws.Cells.Item(2, 13).Value = -738
$ws.Cells.Item(16, 8).Value = 730.3
$ws.Cells.Item(16, 9).Value = 730.3
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 730.3
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -560.3
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 4319.1665
$ws.Cells.Item(31, 10).Value = 3633.3333
$ws.Cells.Item(31, 12).Value = 3633.3333
$ws.Cells.Item(31, 14).Value = -4129.3333
$ws.Cells.Item(40, 8).Value = 4451.8184
$ws.Cells.Item(40, 9).Value = 4333.737
$ws.Cells.Item(40, 11).Value = 4333.737
$ws.Cells.Item(40, 13).Value = -4197.737
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(63, 8).Value = 38000
$ws.Cells.Item(63, 9).Value = 38000
$ws.Cells.Item(63, 11).Value = 38000
$ws.Cells.Item(63, 13).Value = -37251
$ws.Cells.Item(66, 8).Value = 38000
$ws.Cells.Item(66, 9).Value = 38000
$ws.Cells.Item(66, 11).Value = 114000
$ws.Cells.Item(66, 13).Value = -110256
$ws.Cells.Item(82, 8).Value = 1053.5454
$ws.Cells.Item(82, 9).Value = 783.625
$ws.Cells.Item(82, 11).Value = 783.625
$ws.Cells.Item(82, 13).Value = -422.625
$ws.Cells.Item(85, 8).Value = 1053.5454
$ws.Cells.Item(85, 9).Value = 783.625
$ws.Cells.Item(85, 11).Value = 783.625
$ws.Cells.Item(85, 13).Value = 464.375
$ws.Cells.Item(100, 8).Value = 4754.7144
$ws.Cells.Item(100, 9).Value = 3857.2
$ws.Cells.Item(100, 11).Value = 3857.2
$ws.Cells.Item(100, 13).Value = -3316.2
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(115, 8).Value = 50000
$ws.Cells.Item(115, 10).Value = 50000
$ws.Cells.Item(115, 12).Value = 50000
$ws.Cells.Item(115, 14).Value = -52350
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 5716.5557
$ws.Cells.Item(122, 9).Value = 5494.125
$ws.Cells.Item(122, 11).Value = 16482.375
$ws.Cells.Item(122, 13).Value = -14032.375
$ws.Cells.Item(123, 8).Value = 20000
$ws.Cells.Item(123, 9).Value = 20000
$ws.Cells.Item(123, 11).Value = 20000
$ws.Cells.Item(123, 13).Value = -15100
$ws.Cells.Item(130, 9).Value = 100000
$ws.Cells.Item(130, 11).Value = 100000
$ws.Cells.Item(130, 13).Value = -94980
$ws.Cells.Item(132, 8).Value = 3590.9473
$ws.Cells.Item(132, 9).Value = 3648.7334
$ws.Cells.Item(132, 11).Value = 10946.2002
$ws.Cells.Item(132, 13).Value = -8416.200199999999
$ws.Cells.Item(136, 8).Value = 2834.4375
$ws.Cells.Item(136, 9).Value = 2623.4
$ws.Cells.Item(136, 11).Value = 7870.200000000001
$ws.Cells.Item(136, 13).Value = -5320.200000000001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 48000
$ws.Cells.Item(16, 10).Value = 48000
$ws.Cells.Item(16, 12).Value = 48000
$ws.Cells.Item(16, 14).Value = -48584
$ws.Cells.Item(55, 8).Value = 6499
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 4919.8
$ws.Cells.Item(62, 9).Value = 4199.6665
$ws.Cells.Item(62, 10).Value = 6000
$ws.Cells.Item(62, 11).Value = 4199.6665
$ws.Cells.Item(62, 12).Value = 6000
$ws.Cells.Item(62, 13).Value = -3575.6665
$ws.Cells.Item(62, 14).Value = -7248
$ws.Cells.Item(65, 8).Value = 4919.8
$ws.Cells.Item(65, 9).Value = 4199.6665
$ws.Cells.Item(65, 10).Value = 6000
$ws.Cells.Item(65, 11).Value = 20998.3325
$ws.Cells.Item(65, 12).Value = 30000
$ws.Cells.Item(65, 13).Value = -17878.3325
$ws.Cells.Item(65, 14).Value = -36240
$ws.Cells.Item(81, 8).Value = 3667.6
$ws.Cells.Item(81, 10).Value = 4994.25
$ws.Cells.Item(81, 12).Value = 9988.5
$ws.Cells.Item(81, 14).Value = -12110.5
$ws.Cells.Item(84, 8).Value = 3667.6
$ws.Cells.Item(84, 10).Value = 4994.25
$ws.Cells.Item(84, 12).Value = 49942.5
$ws.Cells.Item(84, 14).Value = -60550.5
$ws.Cells.Item(96, 8).Value = 2514.1428
$ws.Cells.Item(96, 9).Value = 1899.75
$ws.Cells.Item(96, 11).Value = 1899.75
$ws.Cells.Item(96, 13).Value = -526.75
$ws.Cells.Item(100, 8).Value = 3669027.8
$ws.Cells.Item(100, 9).Value = 5808657.5
$ws.Cells.Item(100, 10).Value = 1091.4286
$ws.Cells.Item(100, 11).Value = 11617315
$ws.Cells.Item(100, 12).Value = 2182.8572
$ws.Cells.Item(100, 13).Value = -11616774
$ws.Cells.Item(100, 14).Value = -3264.8572
$ws.Cells.Item(112, 8).Value = 49888.668
$ws.Cells.Item(112, 10).Value = 49888.668
$ws.Cells.Item(112, 12).Value = 49888.668
$ws.Cells.Item(112, 14).Value = -52842.668
$ws.Cells.Item(121, 8).Value = 21000
$ws.Cells.Item(121, 10).Value = 21000
$ws.Cells.Item(121, 12).Value = 21000
$ws.Cells.Item(121, 14).Value = -24494
$ws.Cells.Item(122, 8).Value = 4623.75
$ws.Cells.Item(122, 9).Value = 4498.3335
$ws.Cells.Item(122, 11).Value = 13495.0005
$ws.Cells.Item(122, 13).Value = -11045.0005
$ws.Cells.Item(132, 8).Value = 4900
$ws.Cells.Item(132, 9).Value = 4800
$ws.Cells.Item(132, 11).Value = 14400
$ws.Cells.Item(132, 13).Value = -11870
$ws.Cells.Item(136, 8).Value = 1877.6538
$ws.Cells.Item(136, 9).Value = 901.4091
$ws.Cells.Item(136, 10).Value = 7247
$ws.Cells.Item(136, 11).Value = 2704.2273
$ws.Cells.Item(136, 12).Value = 21741
$ws.Cells.Item(136, 13).Value = -154.2273
$ws.Cells.Item(136, 14).Value = -26841

